# Remove the hidden "ACCESS" column (column C). This shifts the old
# "FINAL MH" column (D) left into C, dropping the sheet from A1:D1 to A1:C1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Delete()

# Column C (previously D, "FINAL MH") should be visible with width 12.27,
# matching what the old column D used.
$ws.Columns.Item(3).Hidden = $false
$ws.Columns.Item(3).ColumnWidth = 12.27

# Header row height shrinks now that ACCESS (and its wrapped text) is gone.
$ws.Rows.Item(1).RowHeight = 14.9

# Update the view: scrolled one column right, with C18 selected.
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollColumn = 2
